# implemented status column to TransactionData and items for future workbook output
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

$ws.Range("A13").Value = "Status_Success"
$ws.Range("B13").Value = "Success"
$ws.Range("C13").Value = "Status message for TransactionData to record successful transaction."

$ws.Range("A14").Value = "Status_Failure"
$ws.Range("B14").Value = "Failed"
$ws.Range("C14").Value = "Status message for TransactionData to record failed transaction."

$ws.Range("A15").Value = "Status_Pending"
$ws.Range("B15").Value = "Pending"
$ws.Range("C15").Value = "Status message for TransactionData to record pending transaction."

$ws.Activate()
$ws.Range("A13:C15").Select()
